$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(0.25, 0.125, 0.125, 0.125, 0, 0.125, 0, 0.25, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Activate()
$ws.Range("B9").Select()
